# Apply the data fix described in the commit:
# "Fixed issue with string concatenation in excel_utils.py"
#
# The underlying bug caused each form submission to be written one row
# "behind" (row N held the data meant for row N+1, etc.), and duplicated
# the last row's data across several trailing rows. This edit re-aligns
# the rows 2-9 with their correct values and removes the now-stale
# duplicate rows 10-12, shrinking the used range from A1:E12 to A1:E9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected data for rows 2 through 9 (columns A:E).
$data = @(
    @("pramod",     "563456235263", "prajwalsridhar1999@gmail.com", "sgdf",  "Shree photos sbm road maddur,mandya"),
    @("gt1",         "hgsfz",        "ghdsf@g.com",                  "dgsf",  "Shree photos sbm road maddur,mandya"),
    @("Prajwal S",   "07019660148",  "prajwalsridhar1999@gmail.com", "poiuy", "Shree photos sbm road maddur,mandya"),
    @("Prajwal S",   "07019660148",  "prajwalsridhar1999@gmail.com", "6557",  "Shree photos sbm road maddur,mandya"),
    @("rak",         "6345635",      "prajwalsridhar1999@gmail.com", "ann",   "Shree photos sbm road maddur,mandya"),
    @("Prajwal S",   "07019660148",  "prajwalsridhar1999@gmail.com", "fd",    "Shree photos sbm road maddur,mandya"),
    @("Prajwal S22", "07019660148",  "prajwalsridhar1999@gmail.com", "ds",    "Shree photos sbm road maddur,mandya"),
    @("1",           "1",            "prajwalsridhar1999@gmail.com", "1",     "1")
)

$startRow = 2
$endRow = $startRow + $data.Length - 1

# Format as Text first so numeric-looking strings (phone numbers, etc.)
# keep their leading zeros / exact text form instead of being coerced to
# numbers by Excel's usual "smart" Value assignment.
$ws.Range("A$startRow`:E$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}

# Rows 10-12 were stale duplicate rows; remove them entirely so the
# sheet's used range shrinks back down to A1:E9.
$ws.Range("A10:E12").Delete()
